# videogame_movies.xlsx — "sorted order of video game movies"
#
# The data rows (A2:L8) get re-sorted in descending order of column L
# (ROI). The sheet's zoom/selection and the sort itself (which leaves a
# <sortState> behind, like Excel's Data > Sort dialog does) are updated
# to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds release dates. Re-assign them as real date serials
# (same values already on the sheet) before sorting so the date values
# travel correctly with their row and keep their original date
# formatting/style.
$ws.Range("B2").Value = 41166   # 2012-09-14
$ws.Range("B3").Value = 40431   # 2010-09-10
$ws.Range("B4").Value = 43175   # 2018-03-16
$ws.Range("B5").Value = 41712   # 2014-03-14
$ws.Range("B6").Value = 42531   # 2016-06-10
$ws.Range("B7").Value = 42237   # 2015-08-21
$ws.Range("B8").Value = 43203   # 2018-04-13

# Sort the data range A2:L8 by column L (ROI), descending — same as
# Data > Sort in the Excel UI, which records a <sortState>/<sortCondition>
# on the worksheet.
$sortSpec = $ws.Sort
$sortSpec.SortFields.Clear()
$sortSpec.SortFields.Add($ws.Range("L2:L8"), 0, 2, 0, 0)
$sortSpec.SetRange($ws.Range("A2:L8"))
$sortSpec.Header = 0
$sortSpec.Apply()

# View state: zoomed to 132% with M15:M16 selected.
$app = $ws.Application
$win = $app.ActiveWindow
$win.Zoom = 132
$ws.Range("M15:M16").Select()
